$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @{
    2  = "2025-10-17T07:09:28.495611"
    3  = "2025-10-17T07:09:28.495611"
    4  = "2025-10-17T07:09:28.495611"
    5  = "2025-10-17T07:09:28.495611"
    6  = "2025-10-17T07:09:28.495611"
    7  = "2025-10-17T07:09:28.495611"
    8  = "2025-10-17T07:09:28.497421"
    9  = "2025-10-17T07:09:28.497421"
    10 = "2025-10-17T07:09:28.497938"
    11 = "2025-10-17T07:09:28.497938"
    12 = "2025-10-17T07:09:28.497938"
    13 = "2025-10-17T07:09:28.497938"
    14 = "2025-10-17T07:09:28.497938"
    15 = "2025-10-17T07:09:28.497938"
    16 = "2025-10-17T07:09:28.562124"
    17 = "2025-10-17T07:09:28.562124"
    18 = "2025-10-17T07:09:28.563120"
    19 = "2025-10-17T07:09:28.563120"
    20 = "2025-10-17T07:09:28.563120"
    21 = "2025-10-17T07:09:28.563120"
    22 = "2025-10-17T07:09:28.563120"
    23 = "2025-10-17T07:09:28.564120"
    24 = "2025-10-17T07:09:28.564120"
    25 = "2025-10-17T07:09:28.564120"
    26 = "2025-10-17T07:09:28.629598"
    27 = "2025-10-17T07:09:28.629598"
    28 = "2025-10-17T07:09:28.629598"
    29 = "2025-10-17T07:09:28.629598"
    30 = "2025-10-17T07:09:28.629598"
    31 = "2025-10-17T07:09:28.629598"
    32 = "2025-10-17T07:09:28.629598"
    33 = "2025-10-17T07:09:28.629598"
    34 = "2025-10-17T07:09:28.629598"
    35 = "2025-10-17T07:09:28.629598"
    36 = "2025-10-17T07:09:28.629598"
    37 = "2025-10-17T07:09:28.629598"
    38 = "2025-10-17T07:09:28.629598"
    39 = "2025-10-17T07:09:28.629598"
    40 = "2025-10-17T07:09:28.629598"
    41 = "2025-10-17T07:09:28.629598"
    42 = "2025-10-17T07:09:28.629598"
    43 = "2025-10-17T07:09:28.629598"
    44 = "2025-10-17T07:09:28.629598"
    45 = "2025-10-17T07:09:28.629598"
    46 = "2025-10-17T07:09:28.629598"
    47 = "2025-10-17T07:09:28.629598"
    48 = "2025-10-17T07:09:28.629598"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}
